# regen sval data to filter save games
# Update computed statistic columns (B:E, G) for rows 2-12 on Sheet1.
# Column F ("Win") is left untouched - only the filtered statistics changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: RowNumber, B, C, D, E, G
$newValues = @(
    @{ Row = 2;  B = 0.6606524410359556;  C = 1.655778082260271;  D = 3.537761648806719;   E = 0.4942365360607697; G = 6.348428708163715 },
    @{ Row = 3;  B = 0.6606524410359556;  C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697; G = 2.960089034096801 },
    @{ Row = 4;  B = 0.04271373187048222; C = 0.306821227259698;  D = 0.7527432677738641;  E = 0.4942365360607697; G = 1.596514762964814 },
    @{ Row = 5;  B = 1.455362044514542;   C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697; G = 3.754798637575387 },
    @{ Row = 6;  B = 1.455362044514542;   C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697; G = 3.754798637575387 },
    @{ Row = 7;  B = 0.2917716402565462;  C = 0.306821227259698;  D = 0.1494219747398047;  E = 0.4942365360607697; G = 1.242251378316819 },
    @{ Row = 8;  B = 3.286832544864788;   C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697; G = 5.586269137925634 },
    @{ Row = 9;  B = 1.455362044514542;   C = 0.04071648406533734;D = 3.537761648806719;   E = 0.4942365360607697; G = 5.528076713447369 },
    @{ Row = 10; B = 3.286832544864788;   C = 1.655778082260271;  D = 3.537761648806719;   E = 0.4942365360607697; G = 8.974608811992548 },
    @{ Row = 11; B = 1.455362044514542;   C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697; G = 3.754798637575387 },
    @{ Row = 12; B = 3.286832544864788;   C = 10.34677158129881;  D = 261.3203778131603;   E = 10.19245300693656;  G = 285.1464349462605 }
)

foreach ($entry in $newValues) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B   # B - TB
    $ws.Cells.Item($r, 3).Value = $entry.C   # C - d2S
    $ws.Cells.Item($r, 4).Value = $entry.D   # D - K
    $ws.Cells.Item($r, 5).Value = $entry.E   # E - IP
    $ws.Cells.Item($r, 7).Value = $entry.G   # G - sum
}
